$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 = "Tira de conectores macho acodados" (right angle male pin header strip).
# The per-unit cost link/source was re-priced and its link source changed.
#  - C21: cost/unit formula changed from 0.99*G61/5 to 3.17/10
#  - D21 (shared formula =B21*C21) and the downstream subtotal/total (D60/D62)
#    recalculate automatically once C21 changes.
#  - E21: link source label changed from "eBay" to "Aliexpress"
$ws.Range("C21").Formula = "=3.17/10"
$ws.Range("E21").Value = "Aliexpress"

# View state: window scrolled down (top-left visible row -> 10) with G18 selected.
$ws.Range("G18").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
